# Colour-code the "Movies have an ID..." paragraph:
#   - entity names (Movie/Actor/Director) in red
#   - attribute names (ID/title/release date/name/date of birth) in green
# The underlying text is left completely unchanged; only the run
# segmentation and character formatting are modified, matching the
# diff exactly.

$d = $word.ActiveDocument

$RED   = 255        # RGB(255,0,0)   -> w:color w:val="FF0000"
$GREEN = 5287936    # RGB(0,176,80)  -> w:color w:val="00B050"

# Ordered list of (text, color) pairs whose concatenation reproduces the
# original run's text exactly.
$segments = @(
    @{Text="Movie";            Color=$RED},
    @{Text="s have an ";       Color=$null},
    @{Text="ID";                Color=$GREEN},
    @{Text=", a ";              Color=$null},
    @{Text="title ";            Color=$GREEN},
    @{Text="and a ";            Color=$null},
    @{Text="release date";      Color=$GREEN},
    @{Text=". ";                Color=$null},
    @{Text="Actor";             Color=$RED},
    @{Text="s have an ";        Color=$null},
    @{Text="ID";                Color=$GREEN},
    @{Text=", a ";              Color=$null},
    @{Text="name";              Color=$GREEN},
    @{Text=", and a ";          Color=$null},
    @{Text="date of birth";     Color=$GREEN},
    @{Text=". ";                Color=$null},
    @{Text="Director";          Color=$RED},
    @{Text="s have an ";        Color=$null},
    @{Text="ID";                Color=$GREEN},
    @{Text=", a ";              Color=$null},
    @{Text="name ";             Color=$GREEN},
    @{Text="and a ";            Color=$null},
    @{Text="date of birth";     Color=$GREEN},
    @{Text=". Movies have many actors, actors can act in many movies, movies can have more than one "; Color=$null}
)

# Locate the start of the target sentence (the full original run's text).
$target = "Movies have an ID, a title and a release date. Actors have an ID, a name, and a date of birth. Directors have an ID, a name and a date of birth. Movies have many actors, actors can act in many movies, movies can have more than one "

$findRange = $d.Content
$found = $findRange.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the target paragraph text."
}

$base = $findRange.Start
$offset = 0

foreach ($seg in $segments) {
    $len = $seg.Text.Length
    $r = $d.Range($base + $offset, $base + $offset + $len)
    if ($seg.Color) {
        $r.Font.Color = $seg.Color
    }
    $offset = $offset + $len
}
